# "Brownian verification and twin beam dynamics for size 1"
#
# Insert a new results row (row 8) just under the "Sphere / Final Results"
# block, pushing the Cone / Propellor / Custom Sphere blocks down by one
# row each (their internal layout is unchanged, only their row numbers
# shift), and add a new row of data + a new note string "For Hopping
# purposes" for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 8 and everything below it down by one row.
$ws.Rows.Item(8).Insert() | Out-Null

# Populate the newly freed row 8 with the new measurement.
$ws.Range("A8").Value = 15
$ws.Range("B8").Value = 1064
$ws.Range("B8").NumberFormat = "#,##0"
$ws.Range("C8").Value = 11851752494
$ws.Range("C8").NumberFormat = "#,##0"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.1836997782
$ws.Range("G8").Value = "For Hopping purposes"

# Move the sheet's active-cell selection, matching what was left selected
# when the workbook was saved.
$ws.Range("D17").Select() | Out-Null
